$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.613.94'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.25%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.657.21'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.62'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9979'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -2.85%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '47.16'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -4.57%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3272'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -5.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.131'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -7.73%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07053'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -6.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9986'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.983'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -5.50%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '19.54'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -8.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.620'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -6.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.660.68'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.69%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001047'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -7.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06591'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.08%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '78.94'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -6.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.946'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -7.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.79'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -8.97%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.60'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '24.618.23'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.471'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.414'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -13.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '148.04'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.65'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -8.91%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.216'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '125.21'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -5.84%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.068'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.847'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -14.51%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08476'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.74%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.673'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.33'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -11.70%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.279'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.21%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.215'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -7.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06051'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -9.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02234'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -7.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2069'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -7.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.184'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -11.05%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5929'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -8.51%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.852'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.77'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -7.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5619'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -8.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.65'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.94%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.960'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -8.18%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06972'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.81%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.193'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.77%  '
